$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 3.872992333333334
$ws.Range("H2").Value = 11.618977
$ws.Range("I2").Value = 0.03968168904181824
$ws.Range("J2").Value = 0.03968168904181824
$ws.Range("M2").Value = 9.873811666666667
$ws.Range("N2").Value = 29.621435
$ws.Range("O2").Value = 0.01897536961063408
$ws.Range("P2").Value = 0.01897536961063408
$ws.Range("Q2").Value = 38.24119688577723
$ws.Range("R2").Value = 344.170771971995
$ws.Range("S2").Value = 0.0007529747163427492
$ws.Range("T2").Value = 0.000752974716342749

$ws.Range("G3").Value = 3.872992333333334
$ws.Range("H3").Value = 11.618977
$ws.Range("I3").Value = 0.03968168904181824
$ws.Range("J3").Value = 0.03968168904181824
$ws.Range("O3").Value = 0.368560155467396
$ws.Range("P3").Value = 0.368560155467396
$ws.Range("Q3").Value = 742.7618939018066
$ws.Range("R3").Value = 6684.857045116259
$ws.Range("S3").Value = 0.0146250894824614
$ws.Range("T3").Value = 0.01462508948246139

$ws.Range("G4").Value = 3.872992333333334
$ws.Range("H4").Value = 11.618977
$ws.Range("I4").Value = 0.03968168904181824
$ws.Range("J4").Value = 0.03968168904181824
$ws.Range("M4").Value = 140.35703
$ws.Range("N4").Value = 421.07109
$ws.Range("O4").Value = 0.2697364109842271
$ws.Range("P4").Value = 0.2697364109842271
$ws.Range("Q4").Value = 543.6017011194368
$ws.Range("R4").Value = 4892.415310074931
$ws.Range("S4").Value = 0.01070359638393219
$ws.Range("T4").Value = 0.01070359638393218

$ws.Range("G5").Value = 3.872992333333334
$ws.Range("H5").Value = 11.618977
$ws.Range("I5").Value = 0.03968168904181824
$ws.Range("J5").Value = 0.03968168904181824
$ws.Range("M5").Value = 178.3381523333333
$ws.Range("N5").Value = 535.014457
$ws.Range("O5").Value = 0.3427280639377429
$ws.Range("P5").Value = 0.3427280639377429
$ws.Range("Q5").Value = 690.7022967278323
$ws.Range("R5").Value = 6216.320670550489
$ws.Range("S5").Value = 0.01360002845908191
$ws.Range("T5").Value = 0.01360002845908191

$ws.Range("I6").Value = 0.7580514618940433
$ws.Range("J6").Value = 0.7580514618940433
$ws.Range("M6").Value = 9.873811666666667
$ws.Range("N6").Value = 29.621435
$ws.Range("O6").Value = 0.01897536961063408
$ws.Range("P6").Value = 0.01897536961063408
$ws.Range("Q6").Value = 730.5332989553883
$ws.Range("R6").Value = 6574.799690598495
$ws.Range("S6").Value = 0.01438430667332097
$ws.Range("T6").Value = 0.01438430667332096

$ws.Range("I7").Value = 0.7580514618940433
$ws.Range("J7").Value = 0.7580514618940433
$ws.Range("O7").Value = 0.368560155467396
$ws.Range("P7").Value = 0.368560155467396
$ws.Range("S7").Value = 0.2793875646479554
$ws.Range("T7").Value = 0.2793875646479554

$ws.Range("I8").Value = 0.7580514618940433
$ws.Range("J8").Value = 0.7580514618940433
$ws.Range("M8").Value = 140.35703
$ws.Range("N8").Value = 421.07109
$ws.Range("O8").Value = 0.2697364109842271
$ws.Range("P8").Value = 0.2697364109842271
$ws.Range("Q8").Value = 10384.58982397177
$ws.Range("R8").Value = 93461.30841574594
$ws.Range("S8").Value = 0.2044740806726458
$ws.Range("T8").Value = 0.2044740806726458

$ws.Range("I9").Value = 0.7580514618940433
$ws.Range("J9").Value = 0.7580514618940433
$ws.Range("M9").Value = 178.3381523333333
$ws.Range("N9").Value = 535.014457
$ws.Range("O9").Value = 0.3427280639377429
$ws.Range("P9").Value = 0.3427280639377429
$ws.Range("Q9").Value = 13194.69756482209
$ws.Range("R9").Value = 118752.2780833988
$ws.Range("S9").Value = 0.2598055099001211
$ws.Range("T9").Value = 0.2598055099001211

$ws.Range("G10").Value = 18.786417
$ws.Range("H10").Value = 56.359251
$ws.Range("I10").Value = 0.1924808245004516
$ws.Range("J10").Value = 0.1924808245004516
$ws.Range("M10").Value = 9.873811666666667
$ws.Range("N10").Value = 29.621435
$ws.Range("O10").Value = 0.01897536961063408
$ws.Range("P10").Value = 0.01897536961063408
$ws.Range("Q10").Value = 185.493543349465
$ws.Range("R10").Value = 1669.441890145185
$ws.Range("S10").Value = 0.003652394787855661
$ws.Range("T10").Value = 0.00365239478785566

$ws.Range("G11").Value = 18.786417
$ws.Range("H11").Value = 56.359251
$ws.Range("I11").Value = 0.1924808245004516
$ws.Range("J11").Value = 0.1924808245004516
$ws.Range("O11").Value = 0.368560155467396
$ws.Range("P11").Value = 0.368560155467396
$ws.Range("Q11").Value = 3602.856259346006
$ws.Range("R11").Value = 32425.70633411406
$ws.Range("S11").Value = 0.07094076260237901
$ws.Range("T11").Value = 0.070940762602379

$ws.Range("G12").Value = 18.786417
$ws.Range("H12").Value = 56.359251
$ws.Range("I12").Value = 0.1924808245004516
$ws.Range("J12").Value = 0.1924808245004516
$ws.Range("M12").Value = 140.35703
$ws.Range("N12").Value = 421.07109
$ws.Range("O12").Value = 0.2697364109842271
$ws.Range("P12").Value = 0.2697364109842271
$ws.Range("Q12").Value = 2636.80569446151
$ws.Range("R12").Value = 23731.25125015359
$ws.Range("S12").Value = 0.05191908678403669
$ws.Range("T12").Value = 0.05191908678403669

$ws.Range("G13").Value = 18.786417
$ws.Range("H13").Value = 56.359251
$ws.Range("I13").Value = 0.1924808245004516
$ws.Range("J13").Value = 0.1924808245004516
$ws.Range("M13").Value = 178.3381523333333
$ws.Range("N13").Value = 535.014457
$ws.Range("O13").Value = 0.3427280639377429
$ws.Range("P13").Value = 0.3427280639377429
$ws.Range("Q13").Value = 3350.334896743523
$ws.Range("R13").Value = 30153.01407069171
$ws.Range("S13").Value = 0.06596858032618023
$ws.Range("T13").Value = 0.06596858032618022

$ws.Range("G14").Value = 0.9551306666666667
$ws.Range("H14").Value = 2.865392
$ws.Range("I14").Value = 0.009786024563686945
$ws.Range("J14").Value = 0.009786024563686943
$ws.Range("M14").Value = 9.873811666666667
$ws.Range("N14").Value = 29.621435
$ws.Range("O14").Value = 0.01897536961063408
$ws.Range("P14").Value = 0.01897536961063408
$ws.Range("Q14").Value = 9.430780319724445
$ws.Range("R14").Value = 84.87702287751999
$ws.Range("S14").Value = 0.0001856934331147039
$ws.Range("T14").Value = 0.0001856934331147038

$ws.Range("G15").Value = 0.9551306666666667
$ws.Range("H15").Value = 2.865392
$ws.Range("I15").Value = 0.009786024563686945
$ws.Range("J15").Value = 0.009786024563686943
$ws.Range("O15").Value = 0.368560155467396
$ws.Range("P15").Value = 0.368560155467396
$ws.Range("Q15").Value = 183.1748172572409
$ws.Range("R15").Value = 1648.573355315168
$ws.Range("S15").Value = 0.003606738734600217
$ws.Range("T15").Value = 0.003606738734600216

$ws.Range("G16").Value = 0.9551306666666667
$ws.Range("H16").Value = 2.865392
$ws.Range("I16").Value = 0.009786024563686945
$ws.Range("J16").Value = 0.009786024563686943
$ws.Range("M16").Value = 140.35703
$ws.Range("N16").Value = 421.07109
$ws.Range("O16").Value = 0.2697364109842271
$ws.Range("P16").Value = 0.2697364109842271
$ws.Range("Q16").Value = 134.0593036352533
$ws.Range("R16").Value = 1206.53373271728
$ws.Range("S16").Value = 0.002639647143612403
$ws.Range("T16").Value = 0.002639647143612403

$ws.Range("G17").Value = 0.9551306666666667
$ws.Range("H17").Value = 2.865392
$ws.Range("I17").Value = 0.009786024563686945
$ws.Range("J17").Value = 0.009786024563686943
$ws.Range("M17").Value = 178.3381523333333
$ws.Range("N17").Value = 535.014457
$ws.Range("O17").Value = 0.3427280639377429
$ws.Range("P17").Value = 0.3427280639377429
$ws.Range("Q17").Value = 170.3362383302382
$ws.Range("R17").Value = 1533.026144972144
$ws.Range("S17").Value = 0.003353945252359622
$ws.Range("T17").Value = 0.003353945252359621
